# Apply the RS-RR.schema edit: expand the "requestId" row's Description
# and Example cells in the first table.
#
# w:br (manual line break) is represented in the Word text stream as
# Chr(11) (vertical-tab). Inserting that character via Find/Replace
# causes the engine to split the run into multiple <w:t> elements
# joined by <w:br/>.

$d = $word.ActiveDocument
$vtab = [char]11

# --- Description cell -------------------------------------------------
# (Assigning .Text directly -- rather than Find/Replace -- keeps straight
# apostrophes straight instead of smart-quoting them.)
$t = $d.Tables.Item(1)
$descCell = $t.Cell(3, 5)

$newDescription = "Identifiant unique partagé de la demande de ressource,  généré une seule fois par le système du partenaire qui émet la demande " + $vtab + `
    "Il est valorisé comme suit lors de sa création : " + $vtab + `
    "{orgID}.request.{ID unique de la demande dans le système émetteur}" + $vtab + $vtab + `
    "OU - uniquement si un ID unique de la demande n'est pas disponible : " + $vtab + `
    "OrgId émetteur}.request.{senderCaseId}.{numéro d’ordre chronologique}"

$descCell.Range.Text = $newDescription
Write-Host "Description replaced"

# --- Example cell -------------------------------------------------------
$exCell = $t.Cell(3, 6)
$newExample = "fr.health.samu770.request.1249875" + $vtab + "fr.health.samu690.request.DRFR15690242370035.3"
$exCell.Range.Text = $newExample
Write-Host "Example replaced"
